$wb = $excel.ActiveWorkbook

# --- Verizon sheet: add the two test-result rows ---
$ws = $wb.Worksheets.Item("Verizon")
$ws.Range("A1").Value = "Expexted result"
$ws.Range("A2").Value = "Verizon: Accessories"

# Auto-size column A to fit the new text, then select column B
# (mirrors the author clicking the column-B header before saving).
[void]$ws.Columns("A:A").AutoFit()
[void]$ws.Range("B1:B1048576").Select()

# Make Verizon the active sheet/tab (was Expedia before).
[void]$ws.Activate()
